# Apply localized/updated relic data to the __data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("__data")

# Row 6 - Void Orbit
$ws.Range("D6").Value = "虚空环轨矩阵"
$ws.Range("E6").Value = "熵能"
$ws.Range("K6").Value = "轨道体:+3|每秒伤害:+22虚空"
$ws.Range("L6").Value = "fx/relics/void_orbit.png"

# Row 7 - Sigil Halo
$ws.Range("D7").Value = "圣徽光域"
$ws.Range("E7").Value = "辉耀"
$ws.Range("K7").Value = "减速:+35%|持续:+8s"
$ws.Range("L7").Value = "fx/relics/sigil_halo.png"

# Row 8 - Maelstrom Core
$ws.Range("D8").Value = "潮汐漩核"
$ws.Range("E8").Value = "潮汐"
$ws.Range("K8").Value = "拉拽强度:+100|爆裂伤害:+88霜寒"
$ws.Range("L8").Value = "fx/relics/maelstrom.png"

# Row 9 - Seraph Beacon
$ws.Range("D9").Value = "炽天光塔"
$ws.Range("E9").Value = "辉耀"
$ws.Range("K9").Value = "射速:0.6s|伤害:+16光耀"
$ws.Range("L9").Value = "fx/relics/seraph_beacon.png"

# Row 10 - Aegis Bloom
$ws.Range("D10").Value = "护域绽放"
$ws.Range("E10").Value = "护盾"
$ws.Range("K10").Value = "护盾:+60|持续:+8s"
$ws.Range("L10").Value = "fx/relics/aegis_bloom.png"
